# Reverse the order of the comma-separated "Recorded By" entries in column G
# for every row where multiple names/emails are listed.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ",\s*"
        $reversed = $parts[($parts.Length - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
